$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 45, shifting existing rows 45-52 down to 46-53.
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 with data (copy of row 46's style where needed).
$ws.Range("A45").Value = 10
$ws.Range("B45").Value = "Vega Modelo de Temuco"
$ws.Range("C45").Value = "La Araucanía"
$ws.Range("D45").NumberFormat = $ws.Range("D46").NumberFormat
$ws.Range("D45").Value = 44504
$ws.Range("E45").Value = 9
$ws.Range("F45").Value = 100112022
$ws.Range("G45").Value = "Arveja Verde"
$ws.Range("H45").Value = "Sin especificar"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 205
$ws.Range("K45").Value = 17000
$ws.Range("L45").Value = 18000
$ws.Range("M45").Value = 17463
$ws.Range("N45").Value = "$/saco 25 kilos"
$ws.Range("O45").Value = "Provincia de Limarí"
$ws.Range("P45").Value = 699
$ws.Range("Q45").Value = 25
$ws.Range("R45").Value = "Hortaliza"
